# WorkAt 2016/08/19 17:35 ・音乐插入功能
# Rename Sheet1 -> T_TRACK_TYPE_MST, Sheet2 -> T_ALBUM_TYPE_MST,
# remove the now-unused Sheet3, and populate the new
# T_ALBUM_TYPE_MST sheet with its header + data rows.

$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

$wsTrack = $wb.Worksheets.Item("Sheet1")
$wsTrack.Name = "T_TRACK_TYPE_MST"

$wsAlbum = $wb.Worksheets.Item("Sheet2")
$wsAlbum.Name = "T_ALBUM_TYPE_MST"

$null = $wb.Worksheets.Item("Sheet3").Delete()

# Header row
$wsAlbum.Range("A1").Value = "ALBUM_TYPE_ID"
$wsAlbum.Range("B1").Value = "ALBUM_TYPE_NAME"
$wsAlbum.Range("C1").Value = "DESCRIPTION"

# Data rows
$wsAlbum.Range("A2").Value = 10
$wsAlbum.Range("B2").Value = "OP"
$wsAlbum.Range("C2").Value = "片头曲"
# Match the author's original rich-text run split (头/曲 in a distinct
# East-Asian font run) so this becomes its own shared-string entry
# instead of collapsing onto the identical-looking header cell text.
$wsAlbum.Range("C2").Characters(2,1).Font.Name = "ＭＳ Ｐゴシック"
$wsAlbum.Range("C2").Characters(3,1).Font.Name = "ＭＳ Ｐゴシック"

$wsAlbum.Range("A3").Value = 11
$wsAlbum.Range("B3").Value = "ED"
$wsAlbum.Range("C3").Value = "片尾曲"

$wsAlbum.Range("A4").Value = 12
$wsAlbum.Range("B4").Value = "IN"
$wsAlbum.Range("C4").Value = "插入曲"

$wsAlbum.Range("A5").Value = 13
$wsAlbum.Range("B5").Value = "IM"
$wsAlbum.Range("C5").Value = "印象曲"

$wsAlbum.Range("A6").Value = 14
$wsAlbum.Range("B6").Value = "CS"
$wsAlbum.Range("C6").Value = "角色曲"

$wsAlbum.Range("A7").Value = 20
$wsAlbum.Range("B7").Value = "OST"
$wsAlbum.Range("C7").Value = "原声音乐"

$wsAlbum.Range("A8").Value = 30
$wsAlbum.Range("B8").Value = "Sound"
$wsAlbum.Range("C8").Value = "短声音集"

$wsAlbum.Range("A9").Value = 40
$wsAlbum.Range("B9").Value = "Drama"
$wsAlbum.Range("C9").Value = "广播剧(角色)"

$wsAlbum.Range("A10").Value = 50
$wsAlbum.Range("B10").Value = "Radio"
$wsAlbum.Range("C10").Value = "广播节目"
